$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.398.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.39%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.198.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.16%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.627"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "69.10"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.96%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.592"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.20%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.09%  "

# Row 11
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0951"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.60%  "

# Row 12
$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.14%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.48%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.104"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.84%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.521.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.27%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.25%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.883"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.74%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.198.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.68%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.298.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.45%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0955"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.09%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.35%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.77%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.66%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.53%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +21.46%  "

# Row 26
$ws.Range("E26").Value = "  +6.23%  "

# Row 27
$ws.Range("E27").Value = "  +0.16%  "

# Row 28
$ws.Range("E28").Value = "  +3.36%  "

# Row 29
$ws.Range("E29").Value = "  -3.01%  "

# Row 30
$ws.Range("E30").Value = "  -2.49%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "170.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.06%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.40%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.122"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.27%  "

# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.23%  "

# Row 35
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.123"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.33%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0733"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.99%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +15.84%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.37%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.15%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0301"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.01%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.89%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.73%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +18.55%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.69%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.51%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.204"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.29%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.96%  "

# Row 48
$ws.Range("E48").Value = "  +0.11%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.24%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.65%  "

# Row 51
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.22%  "

